$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-9
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03),
# keeping the existing cell formatting (style s="1", date format).
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45233
}
